$wb = $excel.ActiveWorkbook

# metrics_metadata sheet: remove the "baseline" values from column I (rows 2-17)
$wsMeta = $wb.Worksheets.Item("metrics_metadata")
$wsMeta.Range("I2:I17").ClearContents()

# study sheet: update the selection (no longer the active tab)
$wsStudy = $wb.Worksheets.Item("study")
$wsStudy.Range("E57").Select()

# metrics_metadata becomes the active sheet/tab, with selection at M19
$wsMeta.Activate()
$wsMeta.Range("M19").Select()
